# Update "Forecast Comparison" sheet: insert a new "Week_Start_Date" column
# (B) before the existing ASIN column, re-format the Week labels to drop
# leading zeros, and mark is_holiday_week as a proper boolean.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# 1) Insert a new column at B - this shifts ASIN..is_holiday_week (old B..I)
#    one column to the right (new C..J) and grows the used range to A1:J17.
$ws.Columns.Item(2).Insert()

# 2) Header for the newly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# 3) Week_Start_Date values (stored as literal text, like the other label
#    columns, not as Excel date serials) for weeks 1-16.
$ws.Range("B2:B17").NumberFormat = "@"
$ws.Range("B2").Value = "2025-01-05"
$ws.Range("B3").Value = "2025-01-12"
$ws.Range("B4").Value = "2025-01-19"
$ws.Range("B5").Value = "2025-01-26"
$ws.Range("B6").Value = "2025-02-02"
$ws.Range("B7").Value = "2025-02-09"
$ws.Range("B8").Value = "2025-02-16"
$ws.Range("B9").Value = "2025-02-23"
$ws.Range("B10").Value = "2025-03-02"
$ws.Range("B11").Value = "2025-03-09"
$ws.Range("B12").Value = "2025-03-16"
$ws.Range("B13").Value = "2025-03-23"
$ws.Range("B14").Value = "2025-03-30"
$ws.Range("B15").Value = "2025-04-06"
$ws.Range("B16").Value = "2025-04-13"
$ws.Range("B17").Value = "2025-04-20"

# 4) Drop the leading zero from the week labels (W01 -> W1 ... W09 -> W9).
#    W10 .. W16 are already written without a leading zero.
$ws.Range("A2").Value = "W1"
$ws.Range("A3").Value = "W2"
$ws.Range("A4").Value = "W3"
$ws.Range("A5").Value = "W4"
$ws.Range("A6").Value = "W5"
$ws.Range("A7").Value = "W6"
$ws.Range("A8").Value = "W7"
$ws.Range("A9").Value = "W8"
$ws.Range("A10").Value = "W9"

# 5) is_holiday_week (now column J after the insert) becomes a real boolean
#    instead of the numeric 0/1 flag it used to be.
for ($r = 2; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value = [bool]($cell.Value() -ne 0)
}
